$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.515.95"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "'1.876.97"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "'313.47"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'0.4797"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").Value = "'0.3779"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "'0.07390"
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'20.76"
$ws.Range("E11").Value = "  +5.79%  "
$ws.Range("D12").Value = "'0.07870"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").Value = "'1.873.85"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "'5.452"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "'6.612"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "'91.19"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.000008992"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").Value = "'1.013"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'14.97"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").Value = "'27.539.99"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").Value = "'5.149"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'153.94"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "'18.59"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").Value = "'2.028"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "'5.005"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").Value = "'0.08926"
$ws.Range("D31").Value = "'3.325"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").Value = "'1.221"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("D33").Value = "'4.611"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "'0.7528"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'2.708"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("D36").Value = "'0.02079"
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'3.008"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").Value = "'0.5379"
$ws.Range("E40").Value = "  +2.96%  "
$ws.Range("D41").Value = "'7.102"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "'0.1528"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "'8.456"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").Value = "'0.4849"
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").Value = "'10.62"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "'1.014"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "'1.667"
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("D48").Value = "'103.10"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Value = "'67.23"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").Value = "'0.06111"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'0.9018"
$ws.Range("E51").Value = "  +1.95%  "
